# Update the workbook per the commit:
# values on "Descriptif_numerique", "Regression" and "Regression_R2" sheets
# were recalculated (e.g. after reworking the saison/calendrier pipeline),
# so the cached numeric results changed slightly.

$wb = $excel.ActiveWorkbook

# --- Descriptif_numerique sheet ---
$wsDesc = $wb.Worksheets.Item("Descriptif_numerique")

$wsDesc.Range("C3").Value = 411.2954
$wsDesc.Range("F3").Value = 546.7523

$wsDesc.Range("C4").Value = 502.7635
$wsDesc.Range("F4").Value = 545.8373

$wsDesc.Range("C7").Value = 200
$wsDesc.Range("F7").Value = 374.8753

$wsDesc.Range("C8").Value = 700
$wsDesc.Range("F8").Value = 841.183

$wsDesc.Range("C9").Value = 3150
$wsDesc.Range("F9").Value = 3487.1382

# --- Regression sheet ---
$wsReg = $wb.Worksheets.Item("Regression")

$wsReg.Range("B2").Value = 2.349171977025976
$wsReg.Range("B3").Value = 0.999168542714064
$wsReg.Range("B4").Value = 1.082928697017485
$wsReg.Range("B5").Value = 1.295112582232564

# --- Regression_R2 sheet ---
$wsR2 = $wb.Worksheets.Item("Regression_R2")

$wsR2.Range("A2").Value = 0.9996965393722347
